$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> values for columns A (date text), D, E, F, G, H
$rows = @{
    3  = @{ A = "28-07-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 }
    4  = @{ A = "01-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 }
    5  = @{ A = "04-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 }
    6  = @{ A = "08-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 }
    7  = @{ A = "11-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    8  = @{ A = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    9  = @{ A = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    10 = @{ A = "22-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 }
    11 = @{ A = "25-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    12 = @{ A = "29-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    13 = @{ A = "01-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 }
    14 = @{ A = "05-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 }
    15 = @{ A = "08-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    16 = @{ A = "12-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    17 = @{ A = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    18 = @{ A = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    19 = @{ A = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    20 = @{ A = "26-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
    21 = @{ A = "29-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    # Force the date column to be treated as plain text so Excel does not
    # reinterpret strings like "01-08-2022" as a date serial number. We
    # temporarily mark the cell as Text, write the literal string, then
    # reset the cell style back to Normal so no stray formatting is left
    # behind on the cell.
    $cellA = $ws.Range("A$r")
    $cellA.NumberFormat = "@"
    $cellA.Value = $vals.A
    $cellA.Style = "Normal"

    $ws.Range("D$r").Value = $vals.D
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
}
